# Auto-generated script applying scheduled market-data refresh values
# to the per-class Leve profit tables (currentAveragePrice / LevePrice / LeveProfit columns).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1654.7142
$ws.Range("J17").Value = 1654.7142
$ws.Range("L17").Value = 4964.142599999999
$ws.Range("N17").Value = -5300.142599999999
# Row 98
$ws.Range("H98").Value = 47620772
$ws.Range("I98").Value = 58825084
$ws.Range("J98").Value = 2448.5
$ws.Range("K98").Value = 58825084
$ws.Range("L98").Value = 2448.5
$ws.Range("M98").Value = -58823586
$ws.Range("N98").Value = -5444.5
# Row 100
$ws.Range("H100").Value = 8975.666999999999
$ws.Range("I100").Value = 14242
$ws.Range("K100").Value = 14242
$ws.Range("M100").Value = -13701
# Row 122
$ws.Range("H122").Value = 47620772
$ws.Range("I122").Value = 58825084
$ws.Range("J122").Value = 2448.5
$ws.Range("K122").Value = 176475252
$ws.Range("L122").Value = 7345.5
$ws.Range("M122").Value = -176472802
$ws.Range("N122").Value = -12245.5
# Row 126
$ws.Range("H126").Value = 63000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 63000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 63000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -72880
# Row 128
$ws.Range("H128").Value = 89995
$ws.Range("J128").Value = 89995
$ws.Range("L128").Value = 89995
$ws.Range("N128").Value = -99955
# Row 130
$ws.Range("H130").Value = 63332.668
$ws.Range("J130").Value = 63332.668
$ws.Range("L130").Value = 63332.668
$ws.Range("N130").Value = -73372.66800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1208.68
$ws.Range("I2").Value = 1179.5
$ws.Range("K2").Value = 1179.5
$ws.Range("M2").Value = -1066.5
# Row 32
$ws.Range("H32").Value = 6427062.5
$ws.Range("I32").Value = 7953017
$ws.Range("K32").Value = 7953017
$ws.Range("M32").Value = -7952730
# Row 52
$ws.Range("H52").Value = 70000
$ws.Range("J52").Value = 70000
$ws.Range("L52").Value = 70000
$ws.Range("N52").Value = -70636
# Row 98
$ws.Range("H98").Value = 39749
$ws.Range("J98").Value = 39749
$ws.Range("L98").Value = 39749
$ws.Range("N98").Value = -45739
# Row 114
$ws.Range("H114").Value = 110000
$ws.Range("J114").Value = 110000
$ws.Range("L114").Value = 110000
$ws.Range("N114").Value = -118678
# Row 116
$ws.Range("H116").Value = 1208.68
$ws.Range("I116").Value = 1179.5
$ws.Range("K116").Value = 1179.5
$ws.Range("M116").Value = 1114.5
# Row 123
$ws.Range("H123").Value = 80913.336
$ws.Range("J123").Value = 80913.336
$ws.Range("L123").Value = 80913.336
$ws.Range("N123").Value = -90713.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1208.68
$ws.Range("I3").Value = 1179.5
$ws.Range("K3").Value = 1179.5
$ws.Range("M3").Value = -1065.5
# Row 26
$ws.Range("H26").Value = 35235.5
$ws.Range("J26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("N26").Value = -50584
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 134
$ws.Range("H134").Value = 317184.84
$ws.Range("I134").Value = 1507.3462
$ws.Range("K134").Value = 4522.0386
$ws.Range("M134").Value = -1987.0386

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 92
$ws.Range("H92").Value = 53333
$ws.Range("J92").Value = 53333
$ws.Range("L92").Value = 53333
$ws.Range("N92").Value = -58325
# Row 116
$ws.Range("H116").Value = 87827.336
$ws.Range("J116").Value = 87827.336
$ws.Range("L116").Value = 87827.336
$ws.Range("N116").Value = -97005.336
# Row 119
$ws.Range("H119").Value = 61870.5
$ws.Range("J119").Value = 61870.5
$ws.Range("L119").Value = 61870.5
$ws.Range("N119").Value = -71546.5
# Row 132
$ws.Range("H132").Value = 6723.3477
$ws.Range("I132").Value = 2532.0557
$ws.Range("K132").Value = 7596.1671
$ws.Range("M132").Value = -5066.1671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 173.46666
$ws.Range("I2").Value = 110.5
$ws.Range("J2").Value = 177.96428
$ws.Range("K2").Value = 663
$ws.Range("L2").Value = 1067.78568
$ws.Range("M2").Value = -550
$ws.Range("N2").Value = -1293.78568
# Row 4
$ws.Range("H4").Value = 23703730
$ws.Range("I4").Value = 11255509
$ws.Range("K4").Value = 33766527
$ws.Range("M4").Value = -33766415
# Row 38
$ws.Range("H38").Value = 122.42857
$ws.Range("I38").Value = 137
$ws.Range("J38").Value = 116.6
$ws.Range("K38").Value = 411
$ws.Range("L38").Value = 349.8
$ws.Range("M38").Value = -64
$ws.Range("N38").Value = -1043.8
# Row 40
$ws.Range("H40").Value = 288.77777
$ws.Range("I40").Value = 342.85715
$ws.Range("K40").Value = 1371.4286
$ws.Range("M40").Value = -1302.4286
# Row 97
$ws.Range("H97").Value = 8929793
$ws.Range("I97").Value = 11905057
$ws.Range("K97").Value = 35715171
$ws.Range("M97").Value = -35714675
# Row 140
$ws.Range("H140").Value = 179739.47
$ws.Range("I140").Value = 179739.47
$ws.Range("K140").Value = 539218.41
$ws.Range("M140").Value = -534038.41

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 34485384
$ws.Range("I132").Value = 43480810
$ws.Range("K132").Value = 130442430
$ws.Range("M132").Value = -130439900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 22763.908
$ws.Range("I7").Value = 17481
$ws.Range("J7").Value = 27166.334
$ws.Range("K7").Value = 17481
$ws.Range("L7").Value = 27166.334
$ws.Range("M7").Value = -17369
$ws.Range("N7").Value = -27390.334
# Row 56
$ws.Range("H56").Value = 39747.5
$ws.Range("I56").Value = 35000
$ws.Range("K56").Value = 35000
$ws.Range("M56").Value = -34309
# Row 64
$ws.Range("H64").Value = 13383
$ws.Range("J64").Value = 13383
$ws.Range("L64").Value = 13383
$ws.Range("N64").Value = -13833
# Row 67
$ws.Range("H67").Value = 13383
$ws.Range("J67").Value = 13383
$ws.Range("L67").Value = 13383
$ws.Range("N67").Value = -14943
# Row 99
$ws.Range("H99").Value = 37497.5
# Row 119
$ws.Range("H119").Value = 99990
$ws.Range("J119").Value = 99990
$ws.Range("L119").Value = 99990
$ws.Range("N119").Value = -109666
# Row 126
$ws.Range("H126").Value = 22763.908
$ws.Range("I126").Value = 17481
$ws.Range("J126").Value = 27166.334
$ws.Range("K126").Value = 52443
$ws.Range("L126").Value = 81499.00199999999
$ws.Range("M126").Value = -49973
$ws.Range("N126").Value = -86439.00199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 932.75
$ws.Range("I113").Value = 783.2727
$ws.Range("J113").Value = 1261.6
$ws.Range("K113").Value = 2349.8181
$ws.Range("L113").Value = 3784.8
$ws.Range("M113").Value = -179.8181
$ws.Range("N113").Value = -8124.799999999999
# Row 132
$ws.Range("H132").Value = 483164.38
$ws.Range("I132").Value = 6350.1333
$ws.Range("K132").Value = 19050.3999
$ws.Range("M132").Value = -16520.3999
# Row 136
$ws.Range("H136").Value = 7247.1924
$ws.Range("I136").Value = 7922.316
$ws.Range("K136").Value = 23766.948
$ws.Range("M136").Value = -21216.948

